$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 21 (shifts existing rows 21.. down by one, carrying formatting)
$ws.Range("A21").EntireRow.Insert()

# Populate the new row 21 with the new parameter
$ws.Range("A21").Value = "Suspension_Rod_Rext (mm)"
$ws.Range("B21").Value = 9.75

# Update the selected cell / view to match the saved state
$ws.Range("B21").Select()
